# Swap the roles/content of the two worksheets:
#  - The sheet that currently holds the hotel data ("hotel_info", position 1)
#    becomes the "review_info" sheet: only the review header row (25 cols).
#  - The sheet that currently holds the review header ("review_info", position 2)
#    becomes the "hotel_info" sheet: the hotel header+data row, with a new
#    "State" column inserted right after "Hotel_Name" (value "Louisiana").

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("hotel_info")    # A1:I1 header, A2:I2 data
$ws2 = $wb.Worksheets.Item("review_info")   # A1:Y1 header only

# 1) Stash sheet2's current header (the review_info header, 25 cols) into a
#    scratch row far below sheet1's live data, so it survives while we
#    rebuild sheet2's content from sheet1's hotel data.
$ws2.Range("A1:Y1").Copy($ws1.Range("A10:Y10"))

# 1b) Clear sheet2's old content entirely now that it is safely stashed, so
#     no leftover cells (e.g. columns K:Y) survive past the rebuild below.
$ws2.Range("A1:Y1").ClearContents()

# 2) Build the new "hotel_info" content on sheet2, inserting a "State" column
#    right after "Hotel_Name" (old column B -> still B; old columns C:I shift to D:J).
#    a) Copy STR + Hotel_Name (cols A:B) straight across.
$ws1.Range("A1:B2").Copy($ws2.Range("A1:B2"))
#    b) New "State" column (C): header + value.
$ws2.Range("C1").Value = "State"
$ws2.Range("C2").Value = "Louisiana"
#    c) Copy the remaining old hotel columns (City.. Total_Reviews_num, C:I) into D:J.
$ws1.Range("C1:I2").Copy($ws2.Range("D1:J2"))

# 3) Clear sheet1's old hotel header+data (no longer needed there).
$ws1.Range("A1:I2").ClearContents()

# 4) Move the stashed review header (step 1) into its final place on sheet1.
$ws1.Range("A10:Y10").Copy($ws1.Range("A1:Y1"))
$ws1.Range("A10:Y10").ClearContents()

# 5) Rename the sheets to reflect their new content (use a scratch name to
#    dodge the transient name collision while swapping).
$ws1.Name = "__tmp_swap__"
$ws2.Name = "hotel_info"
$ws1.Name = "review_info"
